$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (execution time ratio), bold 11pt
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Summary block rows 14-17: labels in column A, formulas + bold 12pt
# vertically-centered font in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$formulas = @{
    "B14" = "=AVERAGE(N2:N11)"
    "B15" = "=AVERAGE(Z2:Z11)"
    "B16" = "=MIN(N2:N11)"
    "B17" = "=MAX(Z2:Z11)"
}
foreach ($addr in @("B14","B15","B16","B17")) {
    $c = $ws.Range($addr)
    $c.Formula = $formulas[$addr]
    $c.Font.Bold = $true
    $c.Font.Size = 12
    $c.VerticalAlignment = -4108  # xlVAlignCenter
}

$ws.Rows("14:17").RowHeight = 15.6

# Mirror the selection left by the original author
$ws.Range("A14:B17").Select()

# Page setup (printer defaults picked up when the file was resaved)
$ws.PageSetup.PaperSize = 9       # xlPaperA4
$ws.PageSetup.Orientation = 1     # xlPortrait
